$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.498.84"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "2.367.77"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'0.674"
$ws.Range("D6").Value = "'239.57"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("D7").Value = "'73.43"
$ws.Range("E7").Value = "  +6.46%  "
$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "  +19.57%  "
$ws.Range("E10").Value = "  +6.77%  "
$ws.Range("D11").Value = "'29.92"
$ws.Range("E11").Value = "  +12.19%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.107"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.717.62"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("D14").Value = "'16.96"
$ws.Range("E14").Value = "  +7.93%  "
$ws.Range("D15").Value = "'6.80"
$ws.Range("E15").Value = "  +8.86%  "
$ws.Range("D16").Value = "'0.903"
$ws.Range("E16").Value = "  +6.89%  "
$ws.Range("D17").Value = "2.367.37"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "44.497.15"
$ws.Range("E18").Value = "  +2.56%  "
$ws.Range("E19").Value = "  +4.89%  "
$ws.Range("D20").Value = "'77.57"
$ws.Range("E20").Value = "  +4.58%  "
$ws.Range("D21").Value = "'6.49"
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("D22").Value = "'255.24"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("B23").Value = "WEMIXToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D23").Value = "'3.80"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'2.52"
$ws.Range("E25").Value = "  +2.97%  "
$ws.Range("D26").Value = "'10.43"
$ws.Range("E26").Value = "  +4.19%  "
$ws.Range("D27").Value = "'2.25"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("D30").Value = "'174.31"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("E32").Value = "  +5.38%  "
$ws.Range("D33").Value = "'0.0746"
$ws.Range("E33").Value = "  +7.97%  "
$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("D35").Value = "'5.21"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  +7.86%  "
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("D38").Value = "'6.51"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  +6.70%  "
$ws.Range("D40").Value = "'20.11"
$ws.Range("E40").Value = "  +10.11%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'8.86"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("E43").Value = "  +3.21%  "
$ws.Range("D44").Value = "'0.0985"
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.49"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.184"
$ws.Range("E47").Value = "  +11.89%  "
$ws.Range("D48").Value = "'98.89"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").Value = "1.444.39"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "2.591.46"
$ws.Range("E51").Value = "  -0.17%  "
